$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "# PDF" section header used to live at row 53. A new parameter
# (v_shift_x_lab) is being added to the "# AXES" section right above it,
# so insert a fresh row at 52 and push everything below it down by one.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row with the new parameter: name, value, help text.
$ws.Cells.Item(52, 1).Value = "v_shift_x_lab"
$ws.Cells.Item(52, 2).Value = 0.9
$ws.Cells.Item(52, 3).Value = "shift x_lab down (+1 means one line lower)"

# Match the author's final view state: scrolled down to the new row and
# that row selected.
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A52").Select()
